# Applies the "revision, added pyrolysis and additional figures" edit:
#  - Insert a new parameter row "chemical_recycling_pyrolysis" (TRUE) right
#    after "chemical_recycling_gasification", pushing subsequent rows down.
#  - Remove the explanation text for "ccs_process_co2".
#  - Change "electricity_availability" value from the string "default" to
#    boolean TRUE.
#  - Change "iam_scenario" value from "SSP2_SPA2_19I_D" to "default" and add
#    an explanation "default or user-defined".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 10 (chemical_recycling_gasification is
# row 9; fossil_routes, currently row 10, will shift to row 11).
$ws.Rows("10:10").Insert()

# Populate the newly inserted row 10 with the new parameter.
$ws.Range("A10").Value = "chemical_recycling_pyrolysis"
$ws.Range("B10").Value = $true

# ccs_process_co2 is now on row 16 after the insertion; drop its explanation.
$ws.Range("C16").ClearContents()

# electricity_availability (row 17): switch from "default" string to TRUE.
$ws.Range("B17").Value = $true

# iam_scenario (row 19): change value and add an explanation.
$ws.Range("B19").Value = "default"
$ws.Range("C19").Value = "default or user-defined"
